$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Extend the balance-sheet header/data rows (8-9), reusing the same
# bold/bordered header style already used by row 1 and row 5 headers.
$ws.Range("A5:J5").Copy()
$ws.Range("A8:J8").PasteSpecial(-4122)

$ws.Range("A8").Value = "Shareholders Equity"
$ws.Range("B8").Value = "Total Assets"
$ws.Range("C8").Value = "Current Assets"
$ws.Range("D8").Value = "Assets Non-Current"
$ws.Range("E8").Value = "Current Liabilities"
$ws.Range("F8").Value = "Liabilities Non-Current"
$ws.Range("G8").Value = "Tax Liabilities"
$ws.Range("H8").Value = "Tax Assets"
$ws.Range("I8").Value = "Cash and Equivalents (USD)"
$ws.Range("J8").Value = "Total Liabilities"

$naValues = @("N/A", "N/A", "N/A", "N/A", "N/A", "N/A", "N/A", "N/A", "N/A", "N/A")
for ($c = 1; $c -le 10; $c++) {
    $ws.Cells.Item(9, $c).Value = $naValues[$c - 1]
}
